$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells that will receive numeric-looking strings,
# so Excel stores them as text (matching the original inlineStr type)
# instead of auto-converting to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("E51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = '36.532.22'
$ws.Range("E2").Value = '  +0.34%  '
$ws.Range("D3").Value = '1.959.54'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("D5").Value = '244.10'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = '58.63'
$ws.Range("E7").Value = '  +2.88%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +4.14%  '
$ws.Range("D10").Value = '0.0788'
$ws.Range("E10").Value = '  -5.83%  '
$ws.Range("E11").Value = '  -0.96%  '
$ws.Range("D12").Value = '14.23'
$ws.Range("E12").Value = '  +6.17%  '
$ws.Range("D13").Value = '0.838'
$ws.Range("E13").Value = '  +4.43%  '
$ws.Range("D14").Value = '2.248.77'
$ws.Range("E14").Value = '  +1.40%  '
$ws.Range("D15").Value = '21.26'
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("E16").Value = '  +2.80%  '
$ws.Range("D17").Value = '1.968.99'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = '36.496.69'
$ws.Range("E18").Value = '  +0.63%  '
$ws.Range("D19").Value = '69.72'
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").Value = '0.0₃0848'
$ws.Range("E20").Value = '  -1.54%  '
$ws.Range("D21").Value = '229.86'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("E22").Value = '  +2.04%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("E24").Value = '  +4.88%  '
$ws.Range("E25").Value = '  +3.44%  '
$ws.Range("D26").Value = '9.15'
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("E27").Value = '  +6.32%  '
$ws.Range("D28").Value = '160.69'
$ws.Range("E28").Value = '  -0.09%  '
$ws.Range("D29").Value = '19.44'
$ws.Range("E29").Value = '  +1.32%  '
$ws.Range("E30").Value = '  +2.30%  '
$ws.Range("D31").Value = '1.20'
$ws.Range("E31").Value = '  +7.98%  '
$ws.Range("E32").Value = '  +3.52%  '
$ws.Range("D33").Value = '0.0612'
$ws.Range("E33").Value = '  -2.35%  '
$ws.Range("D34").Value = '4.41'
$ws.Range("E34").Value = '  +6.03%  '
$ws.Range("E35").Value = '  +17.34%  '
$ws.Range("E36").Value = '  +7.53%  '
$ws.Range("E37").Value = '  -0.11%  '
$ws.Range("D38").Value = '1.76'
$ws.Range("E38").Value = '  -1.32%  '
$ws.Range("D39").Value = '5.44'
$ws.Range("E39").Value = '  -10.36%  '
$ws.Range("D40").Value = '0.0979'
$ws.Range("E40").Value = '  +0.89%  '
$ws.Range("E41").Value = '  +1.24%  '
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("D43").Value = '0.0210'
$ws.Range("E43").Value = '  +1.09%  '
$ws.Range("D44").Value = '1.370.72'
$ws.Range("E44").Value = '  +2.89%  '
$ws.Range("D45").Value = '15.74'
$ws.Range("E45").Value = '  +1.40%  '
$ws.Range("D46").Value = '88.09'
$ws.Range("E46").Value = '  +2.27%  '
$ws.Range("E47").Value = '  +0.98%  '
$ws.Range("D48").Value = '7.12'
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").Value = '2.85'
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("D50").Value = '2.139.27'
$ws.Range("E50").Value = '  +1.45%  '
$ws.Range("D51").Value = '44.02'
$ws.Range("E51").Value = '  -0.03%  '
